$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J: "IP Address2" header + 5 IP values, mirroring the existing
# "IP Address" column (F) layout.
$ws.Range("J2").Value = "IP Address2"
$ws.Range("J3").Value = "13.201.115.209"
$ws.Range("J4").Value = "15.206.116.84"
$ws.Range("J5").Value = "13.201.126.69"
$ws.Range("J6").Value = "13.232.35.243"
$ws.Range("J7").Value = "13.201.67.157"

# New column width for J
$ws.Columns.Item(10).ColumnWidth = 18.6

# Header style: bold font, centered + vertically centered, wrapped text
$ws.Range("J2").Font.Bold = $true
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").VerticalAlignment = -4108
$ws.Range("J2").WrapText = $true

# Update active selection as recorded in the saved workbook
$ws.Range("L4").Select() | Out-Null
